$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 134, shifting the
# existing rows 134:195 down to 137:198.
$ws.Rows("134:136").Insert()

# --- New row 134 ---
$ws.Cells.Item(134, 1).Value = 1
$ws.Cells.Item(134, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(134, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(134, 4).Value = 45176
$ws.Cells.Item(134, 5).Value = 15
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100102
$ws.Cells.Item(134, 8).Value = "Cítricos"
$ws.Cells.Item(134, 9).Value = 100102004
$ws.Cells.Item(134, 10).Value = "Mandarina"
$ws.Cells.Item(134, 11).Value = "Murcott"
$ws.Cells.Item(134, 12).Value = "Primera"
$ws.Cells.Item(134, 13).Value = 200
$ws.Cells.Item(134, 14).Value = 16000
$ws.Cells.Item(134, 15).Value = 17000
$ws.Cells.Item(134, 16).Value = 16400
$ws.Cells.Item(134, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(134, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(134, 19).Value = 820
$ws.Cells.Item(134, 20).Value = 20

# --- New row 135 ---
$ws.Cells.Item(135, 1).Value = 1
$ws.Cells.Item(135, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(135, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(135, 4).Value = 45176
$ws.Cells.Item(135, 5).Value = 15
$ws.Cells.Item(135, 6).Value = "Fruta"
$ws.Cells.Item(135, 7).Value = 100102
$ws.Cells.Item(135, 8).Value = "Cítricos"
$ws.Cells.Item(135, 9).Value = 100102004
$ws.Cells.Item(135, 10).Value = "Mandarina"
$ws.Cells.Item(135, 11).Value = "Murcott"
$ws.Cells.Item(135, 12).Value = "Segunda"
$ws.Cells.Item(135, 13).Value = 180
$ws.Cells.Item(135, 14).Value = 14000
$ws.Cells.Item(135, 15).Value = 15000
$ws.Cells.Item(135, 16).Value = 14556
$ws.Cells.Item(135, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(135, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(135, 19).Value = 728
$ws.Cells.Item(135, 20).Value = 20

# --- New row 136 ---
$ws.Cells.Item(136, 1).Value = 1
$ws.Cells.Item(136, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(136, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(136, 4).Value = 45176
$ws.Cells.Item(136, 5).Value = 15
$ws.Cells.Item(136, 6).Value = "Fruta"
$ws.Cells.Item(136, 7).Value = 100102
$ws.Cells.Item(136, 8).Value = "Cítricos"
$ws.Cells.Item(136, 9).Value = 100102004
$ws.Cells.Item(136, 10).Value = "Mandarina"
$ws.Cells.Item(136, 11).Value = "Murcott"
$ws.Cells.Item(136, 12).Value = "Tercera"
$ws.Cells.Item(136, 13).Value = 60
$ws.Cells.Item(136, 14).Value = 12000
$ws.Cells.Item(136, 15).Value = 13000
$ws.Cells.Item(136, 16).Value = 12167
$ws.Cells.Item(136, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(136, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(136, 19).Value = 608
$ws.Cells.Item(136, 20).Value = 20
